{"js": "// Remove the \" (Changed main)\" runs that were appended after\n// \"This is a Microsoft word document.\" in the first paragraph, restoring\n// the paragraph to just its original single run of text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst first = paragraphs.items[0];\nconst range = first.getRange();\nrange.insertText(\"This is a Microsoft word document.\", \"Replace\");\nawait context.sync();\n", "ps1": "# Remove the \" (Changed main)\" text that was appended after\n# \"This is a Microsoft word document.\" in the first paragraph, restoring\n# the paragraph to just its original text.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \" (Changed main)\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Execute([ref]$find.Text, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$find.Replacement.Text, 2)\n"}
